# Model #1 sequence change
#
# The sheet holds a 6-treatment (A-F) crossover-trial dataset:
#   column A = Treatment received that row (letter A-F)
#   column B = Period (1-6)
#   column C = Plaque Level (measurement)
#   column D = Sequence (the 1-6 randomisation sequence a Subject follows)
#   column E = Subject
#
# Each Sequence group cycles the six treatments across the six periods
# following a Williams-design rotation. This edit corrects the Treatment
# letter recorded for every (Period, Sequence) combination to the proper
# Williams-design rotation, i.e. for 1-based Period p and Sequence s:
#
#   letter index = ( base[p-1] + (s-1) ) mod 6,   base = [0,1,5,2,4,3]
#   letter       = "A".."F"[letter index]
#
# which only changes the rows whose previously recorded Treatment letter
# disagreed with that rotation - matching the data correction captured in
# the commit "Model #1 sequence change". Period/Sequence/Subject/Plaque
# Level values themselves are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$letters = @("A", "B", "C", "D", "E", "F")
$base = @(0, 1, 5, 2, 4, 3)

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Row 1 is the header ("Treatment", "Period", "Plaque Level", "Sequence",
# "Subject"), so data starts on row 2.
$firstDataRow = $firstRow + 1

for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $period = [int]$ws.Cells.Item($r, 2).Value2
    $seq = [int]$ws.Cells.Item($r, 4).Value2

    $idx = ($base[$period - 1] + ($seq - 1)) % 6
    $letter = $letters[$idx]

    if ($ws.Cells.Item($r, 1).Text -ne $letter) {
        $ws.Cells.Item($r, 1).Value = $letter
    }
}
